$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.690.68"
$ws.Range("E2").Value = "'  +4.15%  "
$ws.Range("D3").Value = "'3.502.49"
$ws.Range("E3").Value = "'  +1.92%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'592.98"
$ws.Range("E5").Value = "'  +3.52%  "
$ws.Range("D6").Value = "'169.22"
$ws.Range("E6").Value = "'  +6.07%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("D8").Value = "'3.499.55"
$ws.Range("E8").Value = "'  +1.80%  "
$ws.Range("D9").Value = "'0.571"
$ws.Range("E9").Value = "'  -0.03%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "'  +0.40%  "
$ws.Range("E11").Value = "'  +4.84%  "
$ws.Range("E12").Value = "'  +2.90%  "
$ws.Range("D13").Value = "'4.114.72"
$ws.Range("E13").Value = "'  +2.17%  "
$ws.Range("E14").Value = "'  +0.12%  "
$ws.Range("D15").Value = "'28.13"
$ws.Range("E15").Value = "'  +3.24%  "
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("E16").Value = "'  +2.30%  "
$ws.Range("D17").Value = "'66.699.98"
$ws.Range("E17").Value = "'  +4.09%  "
$ws.Range("D18").Value = "'3.504.21"
$ws.Range("E18").Value = "'  +2.04%  "
$ws.Range("E19").Value = "'  +3.66%  "
$ws.Range("D20").Value = "'14.04"
$ws.Range("E20").Value = "'  +3.03%  "
$ws.Range("D21").Value = "'390.33"
$ws.Range("E21").Value = "'  +1.64%  "
$ws.Range("D22").Value = "'7.98"
$ws.Range("E22").Value = "'  +1.71%  "
$ws.Range("D23").Value = "'73.43"
$ws.Range("E23").Value = "'  +2.64%  "
$ws.Range("E24").Value = "'  +8.29%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "'  -0.15%  "
$ws.Range("E26").Value = "'  +3.23%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("E27").Value = "'  +4.19%  "
$ws.Range("E28").Value = "'  +1.15%  "
$ws.Range("E29").Value = "'  +0.00%  "
$ws.Range("E30").Value = "'  +6.16%  "
$ws.Range("E31").Value = "'  +4.54%  "
$ws.Range("D33").Value = "'23.60"
$ws.Range("E33").Value = "'  +2.28%  "
$ws.Range("D34").Value = "'7.44"
$ws.Range("E34").Value = "'  +6.53%  "
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E36").Value = "'  +5.89%  "
$ws.Range("D37").Value = "'161.15"
$ws.Range("E37").Value = "'  +0.16%  "
$ws.Range("D38").Value = "'0.890"
$ws.Range("E38").Value = "'  +4.63%  "
$ws.Range("E39").Value = "'  +3.22%  "
$ws.Range("D40").Value = "'0.0747"
$ws.Range("E40").Value = "'  +2.77%  "
$ws.Range("D41").Value = "'26.67"
$ws.Range("E41").Value = "'  +2.20%  "
$ws.Range("E42").Value = "'  +5.72%  "
$ws.Range("D43").Value = "'2.836.77"
$ws.Range("E43").Value = "'  +0.25%  "
$ws.Range("D44").Value = "'6.65"
$ws.Range("E44").Value = "'  +3.13%  "
$ws.Range("D45").Value = "'43.44"
$ws.Range("E45").Value = "'  +0.88%  "
$ws.Range("D46").Value = "'26.37"
$ws.Range("E46").Value = "'  -2.29%  "
$ws.Range("D47").Value = "'0.0314"
$ws.Range("E47").Value = "'  +3.07%  "
$ws.Range("E48").Value = "'  +3.44%  "
$ws.Range("D49").Value = "'355.40"
$ws.Range("E49").Value = "'  +6.04%  "
$ws.Range("E50").Value = "'  +2.53%  "
$ws.Range("D51").Value = "'34.23"
$ws.Range("E51").Value = "'  +14.73%  "
